$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.007900000000000001, -0.007900000000000001, -0.0334, 0.1203, -0.09520000000000001),
    @(0.0175, 0.0175, 0.0172, 0.1342, -0.1202),
    @(0.0247, 0.0247, 0.0273, 0.1275, -0.074),
    @(0.0187, 0.0187, 0.014, 0.168, -0.0649),
    @(0.0559, 0.0559, 0.0585, 0.1627, -0.0144),
    @(0.0525, 0.0525, 0.0557, 0.1017, -0.0305),
    @(-0.0052, -0.0052, -0.006, 0.0535, -0.0669)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
